$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Fix inconsistent contributor name: "Harsh" -> "Harshvardhan J. Pandit"
# (applies to every cell in column K, rows 2-30, that currently shows "Harsh")
for ($i = 2; $i -le 30; $i++) {
    $ws.Range("K$i").Value = "Harshvardhan J. Pandit"
}

# Fix inconsistent formatting: make K4:K30 use the same style as K2:K3
# by copying the formatting (not the value) from K2 onto K4:K30.
$ws.Range("K2").Copy() | Out-Null
$ws.Range("K4:K30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
